$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (border/bold/center style) from N1 into the
# two new header cells O1:P1 before filling in their values.
$ws.Range("N1").Copy() | Out-Null
$ws.Range("O1:P1").PasteSpecial(-4122) | Out-Null

# Header row (t+14, t+15)
$ws.Range("O1").Value = 14
$ws.Range("P1").Value = 15

# Data rows 2-9 for the two new columns O (t+14) and P (t+15)
$ws.Range("O2").Value = -0.8397183152996887
$ws.Range("P2").Value = -0.6167749765532855

$ws.Range("O3").Value = -0.4624083309197978
$ws.Range("P3").Value = -0.3977475780327427

$ws.Range("O4").Value = 0.03281388342494291
$ws.Range("P4").Value = -0.01360049289061191

$ws.Range("O5").Value = 0.4130507219582314
$ws.Range("P5").Value = 0.3734246626774081

$ws.Range("O6").Value = -0.3015136450913565
$ws.Range("P6").Value = -0.2926850671390637

$ws.Range("O7").Value = -0.1572945797229056
$ws.Range("P7").Value = -0.1569362483016216

$ws.Range("O8").Value = -0.4406984997042788
$ws.Range("P8").Value = -0.432875370004439

$ws.Range("O9").Value = 0.004067465981761826
$ws.Range("P9").Value = 0.003405145540027914
